# Login With Invalid Username 2
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataSet")

# This host's Hyperlink object doesn't support true in-place edits: a
# single item's .Delete()/.TextToDisplay assignment either no-ops or
# appends a stray duplicate hyperlink, so the whole collection has to be
# rebuilt via a collection-level Delete + Add, preserving each link's
# target/order/display text. Re-adding a hyperlink also re-stamps the
# built-in "Hyperlink" cell style as a fresh style entry, so immediately
# copy the original Хипервръзка-styled (but link-free) B3 cell's format
# back onto each cell to keep it pointing at the original style.
$fmtSource = $ws.Range("B3")

function Add-Link([string]$cellRef, [string]$address, [string]$display) {
    if ($display) {
        $ws.Hyperlinks.Add($ws.Range($cellRef), $address, "", "", $display) | Out-Null
    } else {
        $ws.Hyperlinks.Add($ws.Range($cellRef), $address) | Out-Null
    }
    $fmtSource.Copy()
    $ws.Range($cellRef).PasteSpecial(-4122) | Out-Null
}

$ws.Hyperlinks.Delete(1) | Out-Null

Add-Link "B2"  "mailto:abv@abv.bg"    ""
Add-Link "B6"  "mailto:Ilian@mail.bg" "Ilian@mail.bg"
Add-Link "B4"  "mailto:abv@abv.bg"    ""
Add-Link "B5"  "mailto:abv@abv.bg"    ""
Add-Link "B8"  "mailto:abv@abv.bg"    ""
# B9's hyperlink keeps its old cached display text ("abv@abv.bg") even
# though the cell's value is about to change to the new invalid username.
Add-Link "B9"  "mailto:abv@abv.bg"    "abv@abv.bg"
Add-Link "B10" "mailto:abv@abv.bg"    ""
Add-Link "B12" "mailto:abv@abv.bg"    ""

# Row 9: switch the test case to an invalid username, with the password
# now stored as text "123" instead of a number.
$ws.Range("B9").Value = "kjdfb"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "123"

# Move the active selection to the single cell D10.
$ws.Range("D10").Select()
